$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the new rows ---
# Insert from the bottom anchor upward so row indices used below refer to
# ORIGINAL (pre-shift) row numbers at the moment each Insert() runs.
#
# Point B (original row 14, i.e. right after "flammy5"): 10 new rows
#   -> 6 rows that will hold new channel data (future rows 17-22)
#   -> 4 blank rows that grow the gap before the "French" section header
$ws.Range("A14:A23").EntireRow.Insert()

# Point A (original row 13, i.e. right before "flammy5"): 3 new rows
#   -> 3 rows that will hold new channel data (future rows 13-15)
$ws.Range("A13:A15").EntireRow.Insert()

# --- Clear the formatting Excel auto-copied onto the freshly inserted rows ---
# (Insert() copies the format of the row above, which for row 13 would be
# the "English" section-header style; our new rows should be plain/unstyled.)
# Only touch columns A:C here - column D should stay untouched except where
# a value is actually written (row 17), so no stray empty <c> cells appear.
$ws.Range("A13:C22").ClearFormats()

# --- Populate the 3 new rows inserted ABOVE "flammy5" (rows 13-15) ---
$ws.Range("A13").Value = "Clash of Clans with Cam"
$ws.Range("B13").Value = "UCT2x1vuvgYdhk-kQdlzn6yA"
$ws.Range("C13").Value = 260942

$ws.Range("A14").Value = "HaVoC Gaming - Clash of Clans"
$ws.Range("B14").Value = "UC99fa54IUf9RVpEHaUTxQ9w"
$ws.Range("C14").Value = 173755

$ws.Range("A15").Value = "Clash Of Clans - ClashOnGan"
$ws.Range("B15").Value = "UC5DOhpvPfaUfMdzkQ-9fb5g"
$ws.Range("C15").Value = 87492

# --- Row 16 is the pre-existing "flammy5" row (shifted down from row 13) ---
# left untouched: A16/B16/C16/D16 already hold their original values.

# --- Populate the 6 new rows inserted BELOW "flammy5" (rows 17-22) ---
$ws.Range("A17").Value = "The Clash Of Clans Vidz | Damien Elledge"
$ws.Range("B17").Value = "UCb87__fTO0TdQcqMBxcNFMQ"
$ws.Range("C17").Value = 68647
$ws.Range("D17").Value = "드물게 올라옴"

$ws.Range("A18").Value = "Clash Of Clans | GameDiceHD"
$ws.Range("B18").Value = "UC7Wq4cRGhc1JEF-vr13VsOA"
$ws.Range("C18").Value = 63041

$ws.Range("A19").Value = "Clash of Clans | Eclihpse"
$ws.Range("B19").Value = "UCLAOdac7WmMXQKhOP-8lmrQ"
$ws.Range("C19").Value = 41647

$ws.Range("A20").Value = "Clash Of Clans | Mastersaint"
$ws.Range("B20").Value = "UC_mR72CQd3RVHtmFhlY3O1Q"
$ws.Range("C20").Value = 29701

$ws.Range("A21").Value = "Clash of Clans - HDCOC | TeamDTB - Clash of Clans Base Designs"
$ws.Range("B21").Value = "UCaQP9S6tXRHvGOxel7-KFjw"
$ws.Range("C21").Value = 21074

$ws.Range("A22").Value = "COC Nepal"
$ws.Range("B22").Value = "UCXZ8Ko7yNe9bEMZ8Wri-08A"
$ws.Range("C22").Value = 4579

# --- Column A widens (no longer best-fit autosized) ---
# Excel's ColumnWidth (character units) round-trips through a pixel-quantized
# conversion for this workbook's Normal font (7px Maximum Digit Width), so
# the stored sheet width is `round(ColumnWidth*7)/7`. 241/7 is the input
# that lands on the stored width closest to the target 35.125.
$ws.Columns.Item(1).ColumnWidth = 241/7

# --- Selection moves to D23 ---
$ws.Range("D23").Select()

"done"
